$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'308.24"
$ws.Range("E2").Value = "'-0.18%"
$ws.Range("D3").Value = "'41.06"
$ws.Range("E3").Value = "'0.60%"
$ws.Range("D4").Value = "'5.208"
$ws.Range("E4").Value = "'1.91%"
$ws.Range("D5").Value = "'0.07684"
$ws.Range("E5").Value = "'0.65%"
$ws.Range("E6").Value = "'1.60%"
$ws.Range("D7").Value = "'0.9152"
$ws.Range("E7").Value = "'1.24%"
$ws.Range("D9").Value = "'0.1212"
$ws.Range("E9").Value = "'9.18%"
$ws.Range("E10").Value = "'2.11%"
$ws.Range("D11").Value = "'0.09169"
$ws.Range("E11").Value = "'-0.19%"
$ws.Range("D12").Value = "'0.04215"
$ws.Range("E12").Value = "'0.40%"
$ws.Range("E13").Value = "'-0.12%"
$ws.Range("D14").Value = "'0.001258"
$ws.Range("E14").Value = "'0.05%"
$ws.Range("D15").Value = "'0.005744"
$ws.Range("E15").Value = "'1.05%"
$ws.Range("E16").Value = "'1,903.01%"
$ws.Range("D17").Value = "'3.342"
$ws.Range("E17").Value = "'-0.28%"
$ws.Range("D18").Value = "'4.298"
$ws.Range("E18").Value = "'1.21%"
$ws.Range("D20").Value = "'7.370"
$ws.Range("E20").Value = "'12.49%"
$ws.Range("D21").Value = "'0.1382"
$ws.Range("E21").Value = "'1.23%"
$ws.Range("E22").Value = "'-4.17%"
$ws.Range("D23").Value = "'0.04014"
$ws.Range("E23").Value = "'-1.30%"
$ws.Range("D24").Value = "'0.001264"
$ws.Range("E24").Value = "'2.77%"
$ws.Range("D25").Value = "'0.004272"
$ws.Range("E25").Value = "'3.91%"
$ws.Range("D26").Value = "'0.0001300"
$ws.Range("E26").Value = "'0.01%"
$ws.Range("D38").Value = "'0.02498"
$ws.Range("E38").Value = "'3.44%"
$ws.Range("D39").Value = "'0.05323"
$ws.Range("E39").Value = "'2.53%"
$ws.Range("D40").Value = "'0.007844"
$ws.Range("E40").Value = "'0.81%"
$ws.Range("D41").Value = "'0.1314"
$ws.Range("E41").Value = "'0.90%"
$ws.Range("D42").Value = "'0.006506"
$ws.Range("E42").Value = "'-7.68%"
$ws.Range("D43").Value = "'0.001853"
$ws.Range("E43").Value = "'-5.01%"
$ws.Range("E44").Value = "'-6.22%"
$ws.Range("D45").Value = "'0.3340"
$ws.Range("E45").Value = "'0.35%"
$ws.Range("D46").Value = "'0.00006713"
$ws.Range("E46").Value = "'-3.21%"
$ws.Range("E47").Value = "'-0.01%"
$ws.Range("D48").Value = "'0.3345"
$ws.Range("E48").Value = "'990.41%"
$ws.Range("D49").Value = "'0.003103"
$ws.Range("D50").Value = "'0.00002102"
$ws.Range("E50").Value = "'-0.01%"
$ws.Range("E51").Value = "'-0.01%"
